$d = $word.ActiveDocument

# --- 1. English "Programa" italic paragraph: insert line breaks before items 2-9 ---
$en_find = "1. Eco-innovation: concepts, determinant factors, barriers, types of Eco-innovative agents, category of Eco-innovations. 2.  Eco-innovation metrics: Andersen metrics, Arundel & Kemp metrics, OECD metrics.3. Introduction to products life-cycle: analytical perspective, production chain analysis, Green Supply Chain Management Practices.4.  Eco-innovation in the industry: chemistry, agro-food, metal mechanics.5.  Case study of Eco-innovation projects in Brazil.6. Methods and tools to support the process of Eco-innovation: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, others.7. Early identification of failure as support to Eco-innovation: problem, scene, resources.8. TRIZ as a response to Eco-innovation: inventive principles, engineering parameters, contradictions matrix.9. Methodological proposal for Eco-innovative solutions in technological categories: definition, measurement, analysis, creation"
$en_repl = "1. Eco-innovation: concepts, determinant factors, barriers, types of Eco-innovative agents, category of Eco-innovations. ^l2.  Eco-innovation metrics: Andersen metrics, Arundel & Kemp metrics, OECD metrics.^l3. Introduction to products life-cycle: analytical perspective, production chain analysis, Green Supply Chain Management Practices.^l4.  Eco-innovation in the industry: chemistry, agro-food, metal mechanics.^l5.  Case study of Eco-innovation projects in Brazil.^l6. Methods and tools to support the process of Eco-innovation: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, others.^l7. Early identification of failure as support to Eco-innovation: problem, scene, resources.^l8. TRIZ as a response to Eco-innovation: inventive principles, engineering parameters, contradictions matrix.^l9. Methodological proposal for Eco-innovative solutions in technological categories: definition, measurement, analysis, creation"

$ok1 = $d.Content.Find.Execute($en_find, $true, $false, $false, $false, $false, $true, 1, $false, $en_repl, 2)
Write-Output ("English Programa replace: " + $ok1)

# --- 2. Portuguese "Método" run inside "Avaliação": insert line breaks before items 2-9 ---
$pt_find = "1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.4. Eco inovação na indústria: química, agro alimentos, metal mecânica.5. Estudo de casos de projetos de eco inovação no Brasil.6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar"
$pt_repl = "1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.^l2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.^l3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.^l4. Eco inovação na indústria: química, agro alimentos, metal mecânica.^l5. Estudo de casos de projetos de eco inovação no Brasil.^l6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.^l7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.^l8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.^l9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar"

$ok2 = $d.Content.Find.Execute($pt_find, $true, $false, $false, $false, $false, $true, 1, $false, $pt_repl, 2)
Write-Output ("Portuguese Metodo replace: " + $ok2)

# --- 3. "Norma de recuperação" run: split formula and explanation with a line break ---
$nf_find = "NF= (N1 + N2)/2Onde: NF = nota final; N = nota"
$nf_repl = "NF= (N1 + N2)/2^lOnde: NF = nota final; N = nota"

$ok3 = $d.Content.Find.Execute($nf_find, $true, $false, $false, $false, $false, $true, 1, $false, $nf_repl, 2)
Write-Output ("Norma de recuperacao replace: " + $ok3)
